$d = $word.ActiveDocument

# 1. Remove the "Note" paragraph style from the first (only) paragraph,
#    reverting it back to the default style.
$p = $d.Paragraphs.Item(1)
$p.Range.set_Style("Default Paragraph Style")

# 2. Rename the custom style "MarginNoteRIght" -> "MarginNoteRight"
$s = $d.Styles.Item("MarginNoteRIght")
$s.NameLocal = "MarginNoteRight"
